$p = $ppt.ActivePresentation

# Slide 10: "Questions?" title slide - set the (empty) Title placeholder text
$s10 = $p.Slides.Item(10)
$s10.Shapes.Item(1).TextFrame.TextRange.Text = "Questions?"

# Slide 6: "APIs/technology used" slide - fill the Content placeholder with two lines
$s6 = $p.Slides.Item(6)
$s6.Shapes.Item(2).TextFrame.TextRange.Text = "Google Drive’s API`rGoogle’s Login API"
